$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC001 row: RUN_FLAG -> yes, TEST_DATA -> MEM004,MEM001;MEM004,MEM001
$ws.Range("A2").Value = "yes"
$ws.Range("D2").Value = "MEM004,MEM001;MEM004,MEM001"

# TC002 row: RUN_FLAG -> no
$ws.Range("A3").Value = "no"

# Update selection to F9
$ws.Range("F9").Select()
